# Updates the "cryptos" worksheet with freshly scraped price/volume figures.
# Column D (Price) values that look like plain decimal numbers need to be
# force-formatted as Text first, otherwise Excel would silently coerce them
# to numeric values (e.g. "1.0000" -> 1) and lose the original formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.895.04'
$ws.Range('E2').Value = '  -1.23%  '
$ws.Range('D3').Value = '1.897.42'
$ws.Range('E3').Value = '  -0.83%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.0000'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.7580'
$ws.Range('E5').Value = '  +2.39%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '240.23'
$ws.Range('E6').Value = '  -1.46%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3048'
$ws.Range('E8').Value = '  -3.03%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '25.35'
$ws.Range('E9').Value = '  -6.68%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.06847'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07993'
$ws.Range('E11').Value = '  +0.31%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.7507'
$ws.Range('E12').Value = '  -3.98%  '
$ws.Range('D13').Value = '1.901.74'
$ws.Range('E13').Value = '  -0.33%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.224'
$ws.Range('E14').Value = '  -1.52%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '91.21'
$ws.Range('E15').Value = '  -0.92%  '
$ws.Range('D16').Value = '29.895.34'
$ws.Range('E16').Value = '  -1.25%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '13.91'
$ws.Range('E17').Value = '  -3.16%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '5.952'
$ws.Range('E18').Value = '  +1.86%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '239.94'
$ws.Range('E19').Value = '  -2.14%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.000007725'
$ws.Range('E20').Value = '  -1.78%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '1.000'
$ws.Range('E21').Value = '  +0.02%  '
$ws.Range('D22').Value = '2.150.84'
$ws.Range('E22').Value = '  -0.09%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '1.000'
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '6.982'
$ws.Range('E24').Value = '  +4.57%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '9.241'
$ws.Range('E25').Value = '  -2.82%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '165.51'
$ws.Range('E26').Value = '  +0.11%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '18.75'
$ws.Range('E27').Value = '  -1.44%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.1305'
$ws.Range('E28').Value = '  +2.46%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.031'
$ws.Range('E29').Value = '  -4.40%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.354'
$ws.Range('E30').Value = '  +0.39%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.519'
$ws.Range('E31').Value = '  -1.81%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.295'
$ws.Range('E32').Value = '  -0.89%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.031'
$ws.Range('E33').Value = '  -1.44%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.05355'
$ws.Range('E34').Value = '  +2.58%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.254'
$ws.Range('E35').Value = '  -4.01%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.7283'
$ws.Range('E36').Value = '  -3.20%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.721'
$ws.Range('E37').Value = '  -1.31%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01924'
$ws.Range('E38').Value = '  -1.41%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.773'
$ws.Range('E39').Value = '  -0.87%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '6.198'
$ws.Range('E40').Value = '  -3.14%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.4414'
$ws.Range('E41').Value = '  -2.15%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '72.31'
$ws.Range('E42').Value = '  -4.90%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.914'
$ws.Range('E43').Value = '  -1.67%  '
$ws.Range('E44').Value = '  +0.09%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.8272'
$ws.Range('E45').Value = '  -0.63%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '101.13'
$ws.Range('E46').Value = '  -0.01%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '7.577'
$ws.Range('E47').Value = '  -2.46%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '9.764'
$ws.Range('E48').Value = '  -1.30%  '
$ws.Range('D49').Value = '2.052.11'
$ws.Range('E49').Value = '  +0.07%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '36.17'
$ws.Range('E50').Value = '  -3.27%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.05955'
$ws.Range('E51').Value = '  -0.74%  '
